$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update Binance conversion figures in the summary text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.16 = 53026.45 pesos`n✅ 53026.45 pesos = 13.1 = 965.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas!N10/O10/N12/O12: updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 76
$ws2.Range("O10").Value = 4030.01
$ws2.Range("N12").Value = 4048.5
$ws2.Range("O12").Value = 73.7
